$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.315.23"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.815.44"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'325.43"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").Value = "'0.3667"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "'44.94"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'0.07660"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").Value = "'1.144"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "'0.9999"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'22.00"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'6.321"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "'7.476"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "1.819.98"
$ws.Range("E16").Value = "  +3.95%  "
$ws.Range("D17").Value = "'95.42"
$ws.Range("E17").Value = "  +8.17%  "
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'0.06439"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'17.43"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "'6.234"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "28.321.76"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").Value = "'11.57"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'2.139"
$ws.Range("E25").Value = "  -8.02%  "
$ws.Range("D26").Value = "'160.14"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").Value = "'20.71"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "2.030.15"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("D29").Value = "'2.273"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "'1.202"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "'6.018"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "'0.09126"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "'3.558"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'13.05"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").Value = "'0.02396"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("D37").Value = "'5.222"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'0.2173"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").Value = "'0.6584"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "'0.06197"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'1.199"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'8.045"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'13.88"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'0.6098"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").Value = "'3.734"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "'125.55"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "'2.018"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").Value = "'1.165"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("E51").Value = "  +1.37%  "

Write-Host "Update complete"
